$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated values for column F (dSF) per repulled data
$updates = @{
    2  = -7
    3  = -4
    5  = 3
    6  = 0
    8  = -4
    13 = -10
    14 = -5
    17 = 3
    19 = 2
    20 = 5
    25 = 4
    29 = 1
    32 = 3
    37 = 4
    40 = -1
    43 = -5
    45 = 3
    47 = 1
    50 = -4
    60 = -1
    65 = 2
    68 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
